# Add a new worksheet named "L6" at the end of the workbook, containing a
# per-team summary of Form / Goals scored / Goals conceded / Total Goals
# (last six matches), mirroring the other "Form" / "Goals scored" /
# "Goals conceded" sheets already present in the workbook.

$wb = $excel.ActiveWorkbook

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "L6"

# Header row
$newSheet.Range("B1").Value = "Form"
$newSheet.Range("C1").Value = "Goals scored"
$newSheet.Range("D1").Value = "Goals conceded"
$newSheet.Range("E1").Value = "Total Goals"

# Row index column (A2:A19) should hold the text "1".."18", matching the
# look of the other sheets in the workbook. Force text storage so these
# are not re-interpreted as numbers.
$newSheet.Range("A2:A19").NumberFormat = "@"

$rows = @(
    @("Ajax,W W W D W W", "Ajax,5 2 1 1 2 4", "Ajax,0 1 0 1 0 0", "Ajax,5 3 1 2 2 4"),
    @("AZ Alkmaar,W W W W L W", "AZ Alkmaar,4 2 1 2 0 3", "AZ Alkmaar,1 0 0 0 2 1", "AZ Alkmaar,5 2 1 2 2 4"),
    @("Den Haag,L L L D L W", "Den Haag,1 0 1 0 0 3", "Den Haag,2 5 4 0 3 2", "Den Haag,3 5 5 0 3 5"),
    @("FC Emmen,D D W W W L", "FC Emmen,1 1 3 3 3 0", "FC Emmen,1 1 1 1 1 4", "FC Emmen,2 2 4 4 4 4"),
    @("Feyenoord,D D W W D L", "Feyenoord,1 1 2 2 0 2", "Feyenoord,1 1 0 1 0 3", "Feyenoord,2 2 2 3 0 5"),
    @("For Sittard,W L L L W W", "For Sittard,3 0 0 1 3 3", "For Sittard,1 1 2 3 0 0", "For Sittard,4 1 2 4 3 3"),
    @("Groningen,D L W L L L", "Groningen,1 1 1 0 0 1", "Groningen,1 3 0 2 1 2", "Groningen,2 4 1 2 1 3"),
    @("Heerenveen,L D L W L D", "Heerenveen,1 0 1 2 0 2", "Heerenveen,3 0 2 0 2 2", "Heerenveen,4 0 3 2 2 4"),
    @("Heracles,W D L W L W", "Heracles,2 1 0 4 1 4", "Heracles,1 1 3 0 3 0", "Heracles,3 2 3 4 4 4"),
    @("PSV Eindhoven,D L W W W D", "PSV Eindhoven,1 0 3 2 1 2", "PSV Eindhoven,1 2 0 0 0 2", "PSV Eindhoven,2 2 3 2 1 4"),
    @("Sparta Rotterdam,W D W L W W", "Sparta Rotterdam,2 1 3 0 2 2", "Sparta Rotterdam,0 1 2 2 0 1", "Sparta Rotterdam,2 2 5 2 2 3"),
    @("Twente,L D L L L L", "Twente,1 0 1 0 1 0", "Twente,4 0 2 1 2 3", "Twente,5 0 3 1 3 3"),
    @("Utrecht,W W L D W W", "Utrecht,1 4 1 1 2 3", "Utrecht,0 1 2 1 1 2", "Utrecht,1 5 3 2 3 5"),
    @("Vitesse,W D W D D W", "Vitesse,3 0 2 0 0 2", "Vitesse,1 0 1 0 0 1", "Vitesse,4 0 3 0 0 3"),
    @("VVV Venlo,L L L L L L", "VVV Venlo,1 1 0 0 0 0", "VVV Venlo,3 2 1 2 2 4", "VVV Venlo,4 3 1 2 2 4"),
    @("Waalwijk,L W L L L L", "Waalwijk,0 3 1 0 0 1", "Waalwijk,2 1 3 1 1 3", "Waalwijk,2 4 4 1 1 4"),
    @("Willem II,W D L L W L", "Willem II,3 0 0 0 1 2", "Willem II,1 0 1 4 0 3", "Willem II,4 0 1 4 1 5"),
    @("Zwolle,L W L W W L", "Zwolle,0 2 2 1 2 1", "Zwolle,2 1 3 0 0 2", "Zwolle,2 3 5 1 2 3")
)

# Populate column-by-column (all of "Form" first, then all of "Goals
# scored", etc.) so that new entries land in the shared-string table in
# the same order as the source workbook.
for ($i = 0; $i -lt $rows.Length; $i++) {
    $newSheet.Cells.Item($i + 2, 1).Value = "$($i + 1)"
}
for ($i = 0; $i -lt $rows.Length; $i++) {
    $newSheet.Cells.Item($i + 2, 2).Value = $rows[$i][0]
}
for ($i = 0; $i -lt $rows.Length; $i++) {
    $newSheet.Cells.Item($i + 2, 3).Value = $rows[$i][1]
}
for ($i = 0; $i -lt $rows.Length; $i++) {
    $newSheet.Cells.Item($i + 2, 4).Value = $rows[$i][2]
}
for ($i = 0; $i -lt $rows.Length; $i++) {
    $newSheet.Cells.Item($i + 2, 5).Value = $rows[$i][3]
}

# Move the new sheet to the end of the workbook (Worksheets.Add() inserts
# it before the active sheet by default).
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
